# actus-dictionary.xlsx — "Terms" sheet now all unhidden; column A was the
# "Identifier" column whose header was missing from the frozen title row
# (row 1 was blank in column A while the "Identifier" label itself lived in
# A2, alongside the first term's row). The new dictionary generation run
# promotes "Identifier" into the real header row and shifts every term's
# identifier up one row to line up with the rest of that term's data in
# columns B:AC, which never moved. The last identifier (boundaryMonitoring-
# Cycle, formerly alone in row 125) now lives at the end of row 124, and the
# now-empty row 125 disappears.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
if ($ws.Name -ne "Terms") {
    $ws = $wb.Worksheets.Item("Terms")
    $ws.Activate()
}

# Shift column A (identifiers) up by one row, rows 2:125 -> 1:124, carrying
# values *and* formatting (font/fill/border styles) along with them, just
# like an Excel "Delete, Shift Cells Up" restricted to column A only.
$src = $ws.Range("A2:A125")
$dst = $ws.Range("A1:A124")
$src.Copy($dst)

# The source row that used to hold the lone trailing identifier is now
# completely empty (no other columns were ever populated on row 125) -
# clear both value and formatting so no phantom cell/row remains.
$ws.Range("A125").Clear()

# Return the view to the top-left of the sheet (the previous selection sat
# out at AC1, scrolled to column J).
$ws.Range("A1").Select() | Out-Null
